# Entidade e atributos 0.2
# Add "Nota de frete" and "Nota de entrega" as new top-level ("Entidades")
# list items right after the "Hora Entrega" bullet.

$d = $word.ActiveDocument

# Locate the "Hora Entrega" paragraph robustly (rather than a hard-coded index).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Hora Entrega") {
        $target = $p
    }
}

# Insert "Nota de frete" right after it.
$target.Range.InsertParagraphAfter()
$newP1 = $target.Next()
$newP1.Range.ListFormat.ListLevelNumber = 1
$newP1.Range.Text = "Nota de frete"

# Insert "Nota de entrega" right after that.
$newP1.Range.InsertParagraphAfter()
$newP2 = $newP1.Next()
$newP2.Range.ListFormat.ListLevelNumber = 1
$newP2.Range.Text = "Nota de entrega"
